$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 629.7727
$ws.Range("I17").Value = 453
$ws.Range("J17").Value = 669.05554
$ws.Range("K17").Value = 1359
$ws.Range("L17").Value = 2007.16662
$ws.Range("M17").Value = -1191
$ws.Range("N17").Value = -2343.16662
$ws.Range("H39").Value = 20.09091
$ws.Range("I39").Value = 19.1
$ws.Range("K39").Value = 57.3
$ws.Range("M39").Value = 238.7
$ws.Range("H40").Value = 2326.3333
$ws.Range("I40").Value = 1886.9
$ws.Range("K40").Value = 1886.9
$ws.Range("M40").Value = -1711.9
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("M92").Value = ""
$ws.Range("H106").Value = 90948540
$ws.Range("I106").Value = 100035896
$ws.Range("K106").Value = 100035896
$ws.Range("M106").Value = -100035265
$ws.Range("H137").Value = 1511.4584
$ws.Range("I137").Value = 1374.3
$ws.Range("K137").Value = 4122.9
$ws.Range("M137").Value = -1572.9

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H29").Value = 10900
$ws.Range("J29").Value = 10900
$ws.Range("L29").Value = 10900
$ws.Range("N29").Value = -11516
$ws.Range("H36").Value = 2513
$ws.Range("I36").Value = 2513
$ws.Range("K36").Value = 2513
$ws.Range("M36").Value = -2167
$ws.Range("H37").Value = 21999
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 21999
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 21999
$ws.Range("M37").Value = ""
$ws.Range("N37").Value = -22545
$ws.Range("H63").Value = 1131.125
$ws.Range("I63").Value = 1662.25
$ws.Range("J63").Value = 600
$ws.Range("K63").Value = 1662.25
$ws.Range("L63").Value = 600
$ws.Range("M63").Value = -976.25
$ws.Range("N63").Value = -1972
$ws.Range("H66").Value = 1131.125
$ws.Range("I66").Value = 1662.25
$ws.Range("J66").Value = 600
$ws.Range("K66").Value = 8311.25
$ws.Range("L66").Value = 3000
$ws.Range("M66").Value = -4879.25
$ws.Range("N66").Value = -9864
$ws.Range("H97").Value = 434.73685
$ws.Range("I97").Value = 409
$ws.Range("K97").Value = 409
$ws.Range("M97").Value = 87
$ws.Range("H122").Value = 12004.909
$ws.Range("I122").Value = 7705.4
$ws.Range("K122").Value = 23116.2
$ws.Range("M122").Value = -20666.2
$ws.Range("H135").Value = 130000
$ws.Range("J135").Value = 130000
$ws.Range("L135").Value = 130000
$ws.Range("N135").Value = -140140

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 961
$ws.Range("I36").Value = 961
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 961
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -427
$ws.Range("N36").Value = ""
$ws.Range("H86").Value = 2844.4614
$ws.Range("I86").Value = 2831.5
$ws.Range("K86").Value = 2831.5
$ws.Range("M86").Value = -1708.5
$ws.Range("H89").Value = 2844.4614
$ws.Range("I89").Value = 2831.5
$ws.Range("K89").Value = 14157.5
$ws.Range("M89").Value = -8541.5
$ws.Range("H94").Value = 1093.9032
$ws.Range("I94").Value = 650.5769
$ws.Range("K94").Value = 650.5769
$ws.Range("M94").Value = -199.5769
$ws.Range("H107").Value = 1970.7646
$ws.Range("I107").Value = 1428.1428
$ws.Range("K107").Value = 1428.1428
$ws.Range("M107").Value = 491.8571999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 178.8
$ws.Range("I5").Value = 166.66667
$ws.Range("K5").Value = 166.66667
$ws.Range("M5").Value = -54.66667000000001
$ws.Range("H15").Value = 9714.5
$ws.Range("I15").Value = 14420
$ws.Range("J15").Value = 5009
$ws.Range("K15").Value = 14420
$ws.Range("L15").Value = 5009
$ws.Range("M15").Value = -14250
$ws.Range("N15").Value = -5349
$ws.Range("H31").Value = 5521
$ws.Range("I31").Value = 4352.6
$ws.Range("K31").Value = 4352.6
$ws.Range("M31").Value = -4057.6
$ws.Range("H34").Value = 5521
$ws.Range("I34").Value = 4352.6
$ws.Range("K34").Value = 4352.6
$ws.Range("M34").Value = -4150.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 77197.38
$ws.Range("J2").Value = 467.7143
$ws.Range("L2").Value = 2806.2858
$ws.Range("N2").Value = -3032.2858
$ws.Range("H51").Value = 835.5714
$ws.Range("I51").Value = 958.1667
$ws.Range("J51").Value = 100
$ws.Range("K51").Value = 2874.5001
$ws.Range("L51").Value = 300
$ws.Range("M51").Value = -2414.5001
$ws.Range("N51").Value = -1220
$ws.Range("H75").Value = 542.5
$ws.Range("I75").Value = 185
$ws.Range("K75").Value = 555
$ws.Range("M75").Value = 443
$ws.Range("H78").Value = 542.5
$ws.Range("I78").Value = 185
$ws.Range("K78").Value = 1665
$ws.Range("M78").Value = 3327
$ws.Range("H139").Value = 1417
$ws.Range("I139").Value = 1417
$ws.Range("K139").Value = 4251
$ws.Range("M139").Value = 889

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").Value = ""
$ws.Range("H97").Value = 2155.6924
$ws.Range("I97").Value = 589.7143
$ws.Range("J97").Value = 3982.6667
$ws.Range("K97").Value = 589.7143
$ws.Range("L97").Value = 3982.6667
$ws.Range("M97").Value = -93.71429999999998
$ws.Range("N97").Value = -4974.6667
$ws.Range("H102").Value = 2169.35
$ws.Range("I102").Value = 1229.2667
$ws.Range("K102").Value = 1229.2667
$ws.Range("M102").Value = 392.7333000000001
$ws.Range("H122").Value = 36736.484
$ws.Range("I122").Value = 1953.05
$ws.Range("K122").Value = 5859.15
$ws.Range("M122").Value = -3409.15

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1144.4615
$ws.Range("J16").Value = 1600
$ws.Range("L16").Value = 1600
$ws.Range("N16").Value = -1940
$ws.Range("H22").Value = 1321.75
$ws.Range("I22").Value = 1144.5
$ws.Range("J22").Value = 1499
$ws.Range("K22").Value = 1144.5
$ws.Range("L22").Value = 1499
$ws.Range("M22").Value = -849.5
$ws.Range("N22").Value = -2089
$ws.Range("H27").Value = 1321.75
$ws.Range("I27").Value = 1144.5
$ws.Range("J27").Value = 1499
$ws.Range("K27").Value = 1144.5
$ws.Range("L27").Value = 1499
$ws.Range("M27").Value = -1037.5
$ws.Range("N27").Value = -1713
$ws.Range("H34").Value = 6500
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").Value = ""
$ws.Range("H93").Value = 1065.4
$ws.Range("I93").Value = 1073.25
$ws.Range("K93").Value = 1073.25
$ws.Range("M93").Value = 174.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 19073.666
$ws.Range("J69").Value = 19073.666
$ws.Range("L69").Value = 19073.666
$ws.Range("N69").Value = -20571.666
$ws.Range("H72").Value = 19073.666
$ws.Range("J72").Value = 19073.666
$ws.Range("L72").Value = 57220.99800000001
$ws.Range("N72").Value = -64708.99800000001
$ws.Range("H122").Value = 1543.2273
$ws.Range("I122").Value = 1519.5555
$ws.Range("K122").Value = 4558.666499999999
$ws.Range("M122").Value = -2108.666499999999
$ws.Range("H129").Value = 70130.75
$ws.Range("J129").Value = 58449
$ws.Range("L129").Value = 58449
$ws.Range("N129").Value = -68449

